$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 7.143138311642302)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    4 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    5 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    6 = @(0.0006408296065709695, 0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 0, 0.6468712399813661)
    7 = @(3.286832544864788, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 4.840633575959121)
    8 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
